$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.943.93'
$ws.Range('E2').Value = '  -2.41%  '
$ws.Range('D3').Value = '3.125.44'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'" + '592.49'
$ws.Range('E5').Value = '  -2.83%  '
$ws.Range('D6').Value = "'" + '135.85'
$ws.Range('E6').Value = '  -5.85%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '3.112.71'
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('D9').Value = "'" + '0.517'
$ws.Range('E9').Value = '  -2.83%  '
$ws.Range('D10').Value = "'" + '0.145'
$ws.Range('E10').Value = '  -3.73%  '
$ws.Range('D11').Value = "'" + '5.19'
$ws.Range('E11').Value = '  -3.53%  '
$ws.Range('D12').Value = "'" + '0.456'
$ws.Range('E12').Value = '  -3.97%  '
$ws.Range('D13').Value = "'" + '0.0000246'
$ws.Range('E13').Value = '  -3.75%  '
$ws.Range('D14').Value = "'" + '34.09'
$ws.Range('E14').Value = '  -4.11%  '
$ws.Range('D15').Value = '3.638.99'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').Value = "'" + '0.120'
$ws.Range('E16').Value = '  +2.17%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.132.09'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '62.969.86'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('D19').Value = "'" + '6.69'
$ws.Range('E19').Value = '  -3.02%  '
$ws.Range('D20').Value = "'" + '473.63'
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('D21').Value = "'" + '14.20'
$ws.Range('E21').Value = '  -4.22%  '
$ws.Range('D22').Value = "'" + '0.697'
$ws.Range('E22').Value = '  -3.54%  '
$ws.Range('D23').Value = "'" + '7.69'
$ws.Range('E23').Value = '  -2.01%  '
$ws.Range('D24').Value = "'" + '86.92'
$ws.Range('E24').Value = '  +1.63%  '
$ws.Range('D25').Value = "'" + '12.97'
$ws.Range('E25').Value = '  -5.19%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = "'" + '2.72'
$ws.Range('E27').Value = '  -2.57%  '
$ws.Range('D28').Value = "'" + '7.14'
$ws.Range('E28').Value = '  -3.63%  '
$ws.Range('D29').Value = "'" + '7.89'
$ws.Range('E29').Value = '  -7.42%  '
$ws.Range('D30').Value = "'" + '2.04'
$ws.Range('E30').Value = '  -1.41%  '
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('D32').Value = "'" + '26.64'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').Value = "'" + '0.106'
$ws.Range('E33').Value = '  -9.02%  '
$ws.Range('D34').Value = "'" + '2.52'
$ws.Range('E34').Value = '  -4.70%  '
$ws.Range('E35').Value = '  -2.65%  '
$ws.Range('D36').Value = "'" + '5.81'
$ws.Range('E36').Value = '  -2.77%  '
$ws.Range('D37').Value = "'" + '52.02'
$ws.Range('E37').Value = '  -1.45%  '
$ws.Range('D38').Value = '0.0₃0700'
$ws.Range('E38').Value = '  -6.33%  '
$ws.Range('D39').Value = "'" + '0.0386'
$ws.Range('E39').Value = '  -2.64%  '
$ws.Range('D40').Value = "'" + '420.08'
$ws.Range('E40').Value = '  -7.07%  '
$ws.Range('B41').Value = 'Cosmos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D41').Value = "'" + '8.23'
$ws.Range('E41').Value = '  -1.43%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = "'" + '0.115'
$ws.Range('E42').Value = '  -3.63%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.878.65'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = "'" + '2.66'
$ws.Range('E44').Value = '  -11.62%  '
$ws.Range('D45').Value = "'" + '0.261'
$ws.Range('E45').Value = '  -0.80%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').Value = "'" + '0.998'
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').Value = "'" + '2.13'
$ws.Range('E47').Value = '  -5.45%  '
$ws.Range('D48').Value = "'" + '25.69'
$ws.Range('E48').Value = '  -3.47%  '
$ws.Range('E49').Value = '  -1.37%  '
$ws.Range('D50').Value = "'" + '2.27'
$ws.Range('E50').Value = '  -6.75%  '
$ws.Range('D51').Value = "'" + '119.68'
$ws.Range('E51').Value = '  -1.66%  '
